# Generate Report for Handback
# ------------------------------------------------------------------
# This script mirrors a "handback" report generation pass: the status
# text changes from "Ready for handoff" to "Handed back: in sync with
# en-US" everywhere it appears, and each per-language sheet gets its
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated for the two source docs, with the target
# file column turned into a hyperlink (like the existing source-file
# hyperlinks in column A).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab7eccb140de83cb9813273c17b4bdb41160b377/e2e/"
$doc1 = "8ed99678-a6d1-482d-af05-11047ba8aaca.md"
$doc2 = "9f945c4e-42af-42cd-90b2-35194e932408.md"

# ------------------------------------------------------------------
# Overview sheet: update the per-language status columns (E, F) for
# both rows.
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.08
$wsOverview.Columns.Item(6).ColumnWidth = 29.08

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Row 2 - 8ed99678...
$wsZh.Range("I2").Value = $doc1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($urlBase + $doc1), "", "", $doc1)
$wsZh.Range("J2").Value = "8ed99678-a6d1-482d-af05-11047ba8aaca.376a841a6698dd22d2d7140700db44e476c76f12.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-16 00:43:36"

# Row 3 - 9f945c4e...
$wsZh.Range("I3").Value = $doc2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($urlBase + $doc2), "", "", $doc2)
$wsZh.Range("J3").Value = "9f945c4e-42af-42cd-90b2-35194e932408.7bc168bcbc317b9f882d6de8ffb4845ac8480211.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-16 00:43:36"

$wsZh.Columns.Item(3).ColumnWidth = 29.08
$wsZh.Columns.Item(9).ColumnWidth = 39.2
$wsZh.Columns.Item(10).ColumnWidth = 39.2

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Row 2 - 8ed99678...
$wsDe.Range("I2").Value = $doc1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($urlBase + $doc1), "", "", $doc1)
$wsDe.Range("J2").Value = "8ed99678-a6d1-482d-af05-11047ba8aaca.376a841a6698dd22d2d7140700db44e476c76f12.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-16 00:43:43"

# Row 3 - 9f945c4e...
$wsDe.Range("I3").Value = $doc2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($urlBase + $doc2), "", "", $doc2)
$wsDe.Range("J3").Value = "9f945c4e-42af-42cd-90b2-35194e932408.7bc168bcbc317b9f882d6de8ffb4845ac8480211.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-16 00:43:43"

$wsDe.Columns.Item(3).ColumnWidth = 29.08
$wsDe.Columns.Item(9).ColumnWidth = 39.2
$wsDe.Columns.Item(10).ColumnWidth = 39.2

Write-Output "Handback report generated"
